{"js": "// Update the Data Scientist (HelloFresh) bullet points in the Experience\n// section. The four bullets in that table cell change from:\n//   1. \"Supported demand Forecasting efforts during COVID-19 crisis.\"\n//   2. \"Built predictive models to guide new product development.\"\n//   3. \"Implemented Machine Learning models to improve customer retention.\"\n//   4. \"Experience working with data to track customer behaviour and\n//       improve product performance.\"\n// to:\n//   1. \"Support demand Forecasting efforts during COVID-19 crisis.\"\n//   2. \"Lead the development of Customer360 project.\"\n//   3. \"Build predictive models to guide new product development.\"\n//   4. \"Implement Machine Learning models to track market share of meal\n//       kit companies\"\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The Experience table is the second table in the document (index 1).\nconst experienceTable = tables.items[1];\nconst rows = experienceTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Locate the row holding the HelloFresh \"Data Scientist\" bullet points by\n// reading the text of its (merged) cell instead of hard-coding a row\n// index, so the script is resilient to minor structural drift.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nlet targetCell = null;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    if (cell.body.text.indexOf(\"Supported demand Forecasting\") !== -1) {\n      targetCell = cell;\n      break;\n    }\n  }\n  if (targetCell) break;\n}\n\nif (!targetCell) {\n  throw new Error(\"Could not locate the HelloFresh Data Scientist bullet list.\");\n}\n\nconst paragraphs = targetCell.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nconst replacements = [\n  {\n    match: \"Supported demand Forecasting efforts during COVID-19 crisis.\",\n    text: \"Support demand Forecasting efforts during COVID-19 crisis.\"\n  },\n  {\n    match: \"Built predictive models to guide new product development.\",\n    text: \"Lead the development of Customer360 project.\"\n  },\n  {\n    match: \"Implemented Machine Learning models to improve customer retention.\",\n    text: \"Build predictive models to guide new product development.\"\n  },\n  {\n    match: \"Experience working with data to track customer behaviour and improve product performance.\",\n    text: \"Implement Machine Learning models to track market share of meal kit companies\"\n  }\n];\n\nfor (const p of paragraphs.items) {\n  const currentText = p.text.trim();\n  const replacement = replacements.find((r) => currentText === r.match);\n  if (replacement) {\n    p.insertText(replacement.text, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the Data Scientist (HelloFresh) bullet points in the Experience\n# section. The four bullets in that table cell change from:\n#   1. \"Supported demand Forecasting efforts during COVID-19 crisis.\"\n#   2. \"Built predictive models to guide new product development.\"\n#   3. \"Implemented Machine Learning models to improve customer retention.\"\n#   4. \"Experience working with data to track customer behaviour and\n#       improve product performance.\"\n# to:\n#   1. \"Support demand Forecasting efforts during COVID-19 crisis.\"\n#   2. \"Lead the development of Customer360 project.\"\n#   3. \"Build predictive models to guide new product development.\"\n#   4. \"Implement Machine Learning models to track market share of meal\n#       kit companies\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute(\n        $findText,   # FindText\n        $true,       # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        0,           # Wrap (wdFindStop)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        2            # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\nReplace-ExactText \"Supported demand Forecasting efforts during COVID-19 crisis.\" \"Support demand Forecasting efforts during COVID-19 crisis.\"\nReplace-ExactText \"Built predictive models to guide new product development.\" \"Lead the development of Customer360 project.\"\nReplace-ExactText \"Implemented Machine Learning models to improve customer retention.\" \"Build predictive models to guide new product development.\"\nReplace-ExactText \"Experience working with data to track customer behaviour and improve product performance.\" \"Implement Machine Learning models to track market share of meal kit companies\"\n"}
